$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New activity rows appended to the weekly scoreboard log (rows 364-369,
# continuing week 9 and starting week 10).
$newRows = @(
    @{ A = "Matt";   B = 45515; C = "Ride";    D = 65; E = 13.68; F = 344; G = 11; H = 31; I = 17; J = 3; K = 0; L = "Agile Antelope"; M = 9 },
    @{ A = "Eric";   B = 45515; C = "Workout"; D = 92; E = 0;     F = 0;   G = 38; H = 48; I = 5;  J = 1; K = 0; L = "Wily Hyena";     M = 9 },
    @{ A = "Steven"; B = 45515; C = "Walk";    D = 32; E = 1.58;  F = 161; G = 32; H = 0;  I = 0;  J = 0; K = 0; L = "Brave Leopard";  M = 9 },
    @{ A = "Steven"; B = 45515; C = "Walk";    D = 13; E = 0.77;  F = 36;  G = 13; H = 0;  I = 0;  J = 0; K = 0; L = "Brave Leopard";  M = 9 },
    @{ A = "Steven"; B = 45516; C = "Walk";    D = 23; E = 1.2;   F = 108; G = 23; H = 0;  I = 0;  J = 0; K = 0; L = "Brave Leopard";  M = 10 },
    @{ A = "Matt";   B = 45516; C = "Walk";    D = 10; E = 0.48;  F = 0;   G = 10; H = 0;  I = 0;  J = 0; K = 0; L = "Agile Antelope"; M = 10 }
)

$xlUp = -4162
$lastExistingRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row
$firstNewRow = $lastExistingRow + 1
$r = $firstNewRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $r++
}
$lastNewRow = $r - 1

# Match the date formatting already used in column B (numFmt from the row
# right above the new block) instead of creating a brand-new number format.
$ws.Range("B$lastExistingRow").Copy()
$ws.Range("B${firstNewRow}:B${lastNewRow}").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Move the view down to the newly entered rows, landing on the last cell
# typed (bottom-right of the pasted block), matching a natural data-entry
# scroll position.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 348
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M$lastNewRow").Select()
